# Error Calculations and Plots
# Applies the committed edits to missing_data.xlsx:
#  - Removes the "RM 232" and "SC 92" rows (rows shift up to close the gaps)
#  - Fills in / clears a handful of individual "missing data" cells in columns E/F

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the two rows that were deleted from the sheet ---
# Row 26 is "RM 232"; after it is removed, the old row 28 ("SC 92") becomes row 27.
$ws.Rows.Item(26).Delete()
$ws.Rows.Item(27).Delete()

# --- Newly-imputed / newly-missing values (row numbers below are final, post-delete) ---
$ws.Range("E3").Value = -5.7
$ws.Range("F4").ClearContents()
$ws.Range("E5").ClearContents()
$ws.Range("F9").Value = 17.26
$ws.Range("F10").Value = 16.43
$ws.Range("F13").ClearContents()
$ws.Range("F14").ClearContents()
$ws.Range("E21").Value = -8.699999999999999
$ws.Range("E23").ClearContents()
$ws.Range("E32").Value = -6.4
